# Minor tweaks to existential quantification
#
# Three small text edits across the deck:
#  1. Slide 4  (Existential Proofs): merge the introduction-rule line back
#     into a single run (it was split into 3 runs for "... (e : pred w)").
#  2. Slide 11 (Existential Elimination Inference Rule): reword/relayout
#     the premises line, splitting it into several runs of Courier New text.
#  3. Slide 19 (Proof of Existential Negation (1 of 2)): drop the space
#     right after the "∀" in "∀ (T: Type) ...".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 4 - "Existential Proofs": collapse 3 runs into 1.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)            # "Content Placeholder 2"
$tr4 = $sh4.TextFrame.TextRange

# The introduction-rule line is its own paragraph; find it and replace the
# whole paragraph's characters with the merged text in one go (this makes
# PowerPoint collapse what used to be 3 runs into a single run carrying the
# first run's formatting).
for ($i = 1; $i -le $tr4.Paragraphs().Count; $i++) {
    $para = $tr4.Paragraphs($i, 1)
    if ($para.Text -like "(T : Type) (pred*") {
        $full = $tr4.Characters($para.Start, $para.Length)
        $full.Text = "(T : Type) (pred: T $([char]0x2192) Prop) (w : T) (e : pred w)"
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 11 - "Existential Elimination Inference Rule": reword the line.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)          # "Content Placeholder 2"
$tr11 = $sh11.TextFrame.TextRange

for ($i = 1; $i -le $tr11.Paragraphs().Count; $i++) {
    $para = $tr11.Paragraphs($i, 1)
    if ($para.Text -like "*Q : Prop*P w*") {
        $arrow = [char]0x2192
        $forall = [char]0x2200
        $exists = [char]0x2203

        $newText = "Q: Prop; T: Type; P: (T $arrow Prop); $exists(x: T), P x; ($forall(a: T), P a) $arrow Q"

        # Replace the whole paragraph text first (single run).
        $full = $tr11.Characters($para.Start, $para.Length)
        $full.Text = $newText

        # Re-fetch the (now shorter) paragraph and split it into the 6 runs
        # shown in the target deck by nudging the font size (no value
        # change) at each boundary - this forces PowerPoint to break the
        # run without altering the text itself.
        $para2 = $tr11.Paragraphs($i, 1)
        $base = $para2.Start

        $chunks = @(
            "Q: ",
            "Prop",
            "; T: Type; P: (T $arrow Prop); $exists(x: T), ",
            "P x",
            "; ($forall(a: T), P a) ",
            "$arrow Q"
        )

        $offset = 0
        foreach ($chunk in $chunks) {
            $len = $chunk.Length
            $rng = $tr11.Characters($base + $offset, $len)
            $rng.Font.Size = 18
            $offset += $len
        }
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Slide 19 - "Proof of Existential Negation (1 of 2)": drop one space.
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$sh19 = $s19.Shapes.Item(2)          # "Content Placeholder 2"
$tr19 = $sh19.TextFrame.TextRange

for ($i = 1; $i -le $tr19.Paragraphs().Count; $i++) {
    $para = $tr19.Paragraphs($i, 1)
    if ($para.Text -like "*(T: Type) (pred*") {
        $full = $tr19.Characters($para.Start, $para.Length)
        $full.Text = "  $([char]0x2200)(T: Type) (pred: (T $([char]0x2192) Prop)),"
        break
    }
}
